# Auto-generated edit script: update Leve Profit sheet price snapshots
# per scheduled market-data runner (Mateus_Profits workbook).
$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 3162.923
$ws.Range("I19").Value = 2316.75
$ws.Range("K19").Value = 2316.75
$ws.Range("M19").Value = -2141.75
# row 88
$ws.Range("H88").Value = 4807.75
$ws.Range("J88").Value = 4699.4
$ws.Range("L88").Value = 4699.4
$ws.Range("N88").Value = -5511.4
# row 91
$ws.Range("H91").Value = 4807.75
$ws.Range("J91").Value = 4699.4
$ws.Range("L91").Value = 4699.4
$ws.Range("N91").Value = -7507.4
# row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").ClearContents()
# row 113
$ws.Range("H113").Value = 200004580
$ws.Range("I113").Value = 250002220
$ws.Range("K113").Value = 250002220
$ws.Range("M113").Value = -249998966
# row 132
$ws.Range("H132").Value = 1791.1904
$ws.Range("I132").Value = 1755.75
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5267.25
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2737.25
$ws.Range("N132").Value = -12560
# row 138
$ws.Range("H138").Value = 17243524
$ws.Range("J138").Value = 20835526
$ws.Range("L138").Value = 62506578
$ws.Range("N138").Value = -62516858

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 7384.0527
$ws.Range("I32").Value = 7238.722
$ws.Range("K32").Value = 7238.722
$ws.Range("M32").Value = -6951.722
# row 102
$ws.Range("H102").Value = 3008.7942
$ws.Range("I102").Value = 2196.8
$ws.Range("K102").Value = 2196.8
$ws.Range("M102").Value = -574.8000000000002
# row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# row 135
$ws.Range("H135").Value = 110476.336
$ws.Range("J135").Value = 110476.336
$ws.Range("L135").Value = 110476.336
$ws.Range("N135").Value = -120616.336
# row 139
$ws.Range("H139").Value = 155495
$ws.Range("J139").Value = 155495
$ws.Range("L139").Value = 155495
$ws.Range("N139").Value = -165775

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 3822.9473
$ws.Range("I94").Value = 3109.4666
$ws.Range("K94").Value = 3109.4666
$ws.Range("M94").Value = -2658.4666
# row 134
$ws.Range("H134").Value = 4990.32
$ws.Range("I134").Value = 4781.5835
$ws.Range("K134").Value = 14344.7505
$ws.Range("M134").Value = -11809.7505

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 9089
$ws.Range("I31").Value = 7500
$ws.Range("J31").Value = 10995.8
$ws.Range("K31").Value = 7500
$ws.Range("L31").Value = 10995.8
$ws.Range("M31").Value = -7205
$ws.Range("N31").Value = -11585.8
# row 34
$ws.Range("H34").Value = 9089
$ws.Range("I34").Value = 7500
$ws.Range("J34").Value = 10995.8
$ws.Range("K34").Value = 7500
$ws.Range("L34").Value = 10995.8
$ws.Range("M34").Value = -7298
$ws.Range("N34").Value = -11399.8
# row 58
$ws.Range("H58").Value = 13124.25
$ws.Range("I58").Value = 7500
$ws.Range("J58").Value = 14999
$ws.Range("K58").Value = 7500
$ws.Range("L58").Value = 14999
$ws.Range("M58").Value = -7297
$ws.Range("N58").Value = -15405
# row 136
$ws.Range("H136").Value = 13124.25
$ws.Range("I136").Value = 7500
$ws.Range("J136").Value = 14999
$ws.Range("K136").Value = 22500
$ws.Range("L136").Value = 44997
$ws.Range("M136").Value = -19950
$ws.Range("N136").Value = -50097

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# row 97
$ws.Range("H97").Value = 2875
$ws.Range("I97").Value = 2653.375
$ws.Range("K97").Value = 7960.125
$ws.Range("M97").Value = -7464.125
# row 117
$ws.Range("H117").Value = 85001000
$ws.Range("J117").Value = 68667870
$ws.Range("L117").Value = 206003610
$ws.Range("N117").Value = -206010494
# row 121
$ws.Range("H121").Value = 61905190
$ws.Range("J121").Value = 61905190
$ws.Range("L121").Value = 185715570
$ws.Range("N121").Value = -185718190
# row 131
$ws.Range("H131").Value = 20835302
$ws.Range("J131").Value = 2383.5
$ws.Range("L131").Value = 7150.5
$ws.Range("N131").Value = -17230.5

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 5197.3
$ws.Range("I102").Value = 4784.25
$ws.Range("K102").Value = 4784.25
$ws.Range("M102").Value = -3162.25

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 4359.8
# row 46
$ws.Range("H46").Value = 16501.125
$ws.Range("J46").Value = 20993.25
$ws.Range("L46").Value = 20993.25
$ws.Range("N46").Value = -21369.25
# row 61
$ws.Range("H61").Value = 64634.312
$ws.Range("I61").Value = 73698.14
$ws.Range("J61").Value = 1187.5
$ws.Range("K61").Value = 73698.14
$ws.Range("L61").Value = 1187.5
$ws.Range("M61").Value = -73496.14
$ws.Range("N61").Value = -1591.5
# row 113
$ws.Range("H113").Value = 64634.312
$ws.Range("I113").Value = 73698.14
$ws.Range("J113").Value = 1187.5
$ws.Range("K113").Value = 73698.14
$ws.Range("L113").Value = 1187.5
$ws.Range("M113").Value = -71528.14
$ws.Range("N113").Value = -5527.5
# row 122
$ws.Range("H122").Value = 2922.1365
$ws.Range("I122").Value = 2902.0527
$ws.Range("K122").Value = 8706.158100000001
$ws.Range("M122").Value = -6256.158100000001
# row 132
$ws.Range("H132").Value = 18238.785
$ws.Range("I132").Value = 21885
$ws.Range("J132").Value = 9123.25
$ws.Range("K132").Value = 65655
$ws.Range("L132").Value = 27369.75
$ws.Range("M132").Value = -63125
$ws.Range("N132").Value = -32429.75
# row 136
$ws.Range("H136").Value = 3610.7778
$ws.Range("I136").Value = 3517.64
$ws.Range("K136").Value = 10552.92
$ws.Range("M136").Value = -8002.92

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# row 12
$ws.Range("H12").Value = 13999
$ws.Range("I12").Value = 13999
$ws.Range("K12").Value = 13999
$ws.Range("M12").Value = -13857
# row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
# row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
# row 107
$ws.Range("H107").Value = 1095.6666
$ws.Range("I107").Value = 1015
$ws.Range("K107").Value = 3045
$ws.Range("M107").Value = -1125
# row 132
$ws.Range("H132").Value = 7521.1113
$ws.Range("I132").Value = 7211.25
$ws.Range("K132").Value = 21633.75
$ws.Range("M132").Value = -19103.75

Write-Output "Applied 163 value updates and 7 clears"